# Update BOM: the "1220 Capacitor" device is now specified as
# "1220 Tantalum Capacitor" (row 6, Device column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "1220 Tantalum Capacitor "

# Reflect the active selection recorded in the saved workbook.
$ws.Range("D25").Select() | Out-Null
